# Update SwaadSutra_Daily_2026-01-14.xlsx
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Daily Orders" - add header row + the day's single order row
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Daily Orders")

$headers1 = @("Order ID","Date","Customer","Flat No","Phone","Items","Total","Status","Payment","Collection Date","Collection Time","Notes","Cancel Reason","Feedback")
for ($i = 0; $i -lt $headers1.Count; $i++) {
    $ws1.Cells.Item(1, $i + 1).Value = $headers1[$i]
}

# Plain numbers
$ws1.Cells.Item(2, 1).Value = 14
$ws1.Cells.Item(2, 7).Value = 600

# Text fields. A couple of these (Phone / Collection Date) look like a
# number or an ISO date, so a leading quote keeps Excel from silently
# re-typing them as numeric/date values - they must stay text, exactly
# like the "numberStoredAsText" warning already suppressed on this sheet.
$ws1.Cells.Item(2, 2).Value = "2026-01-14 17:08"
$ws1.Cells.Item(2, 3).Value = "Mrunal"
$ws1.Cells.Item(2, 4).Value = "KLV B 2108"
$ws1.Cells.Item(2, 5).Value = "'9404665203"
$ws1.Cells.Item(2, 6).Value = "Wheat Chapati x40"
$ws1.Cells.Item(2, 8).Value = "NEW"
$ws1.Cells.Item(2, 9).Value = "PENDING"
$ws1.Cells.Item(2, 10).Value = "'2026-01-15"
$ws1.Cells.Item(2, 11).Value = "00:30"
# Notes / Cancel Reason / Feedback are blank for this order

# ---------------------------------------------------------------------
# Sheet "Summary" - bump totals for the new order
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Summary")
$ws2.Cells.Item(2, 1).Value = 1   # Total Orders
$ws2.Cells.Item(2, 2).Value = 1   # New
$ws2.Cells.Item(2, 7).Value = 600 # Total Revenue

# ---------------------------------------------------------------------
# Sheet "Items Breakdown" - add header row + item summary row
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Items Breakdown")

$headers3 = @("Item","Quantity Ordered","Revenue")
for ($i = 0; $i -lt $headers3.Count; $i++) {
    $ws3.Cells.Item(1, $i + 1).Value = $headers3[$i]
}

$ws3.Cells.Item(2, 1).Value = "Wheat Chapati"
$ws3.Cells.Item(2, 2).Value = 40
$ws3.Cells.Item(2, 3).Value = 600
